$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 292
$ws.Range("F4").Value = 665
$ws.Range("F5").Value = 2972
$ws.Range("F7").Value = 244
$ws.Range("F10").Value = 7007
$ws.Range("F12").Value = 113
$ws.Range("F13").Value = 375
$ws.Range("F14").Value = 615
$ws.Range("F15").Value = 1515
$ws.Range("F16").Value = 1132
$ws.Range("F17").Value = 2269
$ws.Range("F18").Value = 1523
$ws.Range("F19").Value = 131
$ws.Range("F20").Value = 1124
$ws.Range("F21").Value = 142
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = 197
$ws.Range("F24").Value = 354
$ws.Range("F25").Value = 21
$ws.Range("F26").Value = 1769
$ws.Range("F27").Value = 1707
$ws.Range("F28").Value = 1039
$ws.Range("F29").Value = 41
$ws.Range("F30").Value = 1680
$ws.Range("F31").Value = 1242
$ws.Range("F32").Value = 146
$ws.Range("F34").Value = 7
$ws.Range("F35").Value = 1068
$ws.Range("F36").Value = 447
$ws.Range("F37").Value = 29
$ws.Range("F38").Value = 2516
$ws.Range("F39").Value = 2761
$ws.Range("F41").Value = 30
$ws.Range("F46").Value = 330
$ws.Range("F48").Value = 175
$ws.Range("F49").Value = 417

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 177
$ws.Range("F10").Value = 33
$ws.Range("F13").Value = 5
$ws.Range("F15").Value = 60
$ws.Range("F19").Value = 48
$ws.Range("F23").Value = 486
$ws.Range("F30").Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1709
$ws.Range("F7").Value = 1858
$ws.Range("F8").Value = 2765
$ws.Range("F9").Value = 1039
$ws.Range("F10").Value = 955
$ws.Range("F13").Value = 1543
$ws.Range("F14").Value = 7432

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 292
$ws.Range("F3").Value = 665
$ws.Range("F4").Value = 2972
$ws.Range("F5").Value = 244
$ws.Range("F6").Value = 1709
$ws.Range("F7").Value = 2765
$ws.Range("F8").Value = 7007
$ws.Range("F9").Value = 1039
$ws.Range("F11").Value = 113
$ws.Range("F12").Value = 375
$ws.Range("F13").Value = 177
$ws.Range("F14").Value = 1543
$ws.Range("F15").Value = 615
$ws.Range("F16").Value = 1515
$ws.Range("F17").Value = 1132
$ws.Range("F18").Value = 2269
$ws.Range("F19").Value = 1523
$ws.Range("F20").Value = 131
$ws.Range("F22").Value = 1124
$ws.Range("F23").Value = 142
$ws.Range("F25").Value = 21
$ws.Range("F26").Value = 1769
$ws.Range("F27").Value = 1039
$ws.Range("F28").Value = 41
$ws.Range("F29").Value = 1680
$ws.Range("F30").Value = 1242
$ws.Range("F31").Value = 146
$ws.Range("F33").Value = 1068
$ws.Range("F36").Value = 486
$ws.Range("F37").Value = 447
$ws.Range("F38").Value = 29
$ws.Range("F39").Value = 2516
$ws.Range("F40").Value = 2762
$ws.Range("F44").Value = 330
$ws.Range("F46").Value = 175
$ws.Range("F47").Value = 417
